# Fruta / hortaliza, semanal
# Update the weekly data rows (D, J, K, L, M, P columns) for the Perejil sheet.
# The underlying change re-shuffles the per-date records across rows 2-19
# (row 12 keeps its original data); this script writes the resulting
# target values directly, cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44260
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(2, 11).Value = 3500
$ws.Cells.Item(2, 12).Value = 3500
$ws.Cells.Item(2, 13).Value = 3500
$ws.Cells.Item(2, 16).Value = 1167

# Row 3
$ws.Cells.Item(3, 4).Value = 44187
$ws.Cells.Item(3, 10).Value = 65
$ws.Cells.Item(3, 11).Value = 3000
$ws.Cells.Item(3, 12).Value = 3000
$ws.Cells.Item(3, 13).Value = 3000
$ws.Cells.Item(3, 16).Value = 1000

# Row 4
$ws.Cells.Item(4, 4).Value = 44223
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 2500
$ws.Cells.Item(4, 12).Value = 3000
$ws.Cells.Item(4, 13).Value = 2781
$ws.Cells.Item(4, 16).Value = 927

# Row 5
$ws.Cells.Item(5, 4).Value = 44222
$ws.Cells.Item(5, 10).Value = 45
$ws.Cells.Item(5, 11).Value = 3000
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 13).Value = 3000
$ws.Cells.Item(5, 16).Value = 1000

# Row 6
$ws.Cells.Item(6, 4).Value = 44389
$ws.Cells.Item(6, 10).Value = 81
$ws.Cells.Item(6, 11).Value = 2800
$ws.Cells.Item(6, 12).Value = 3000
$ws.Cells.Item(6, 13).Value = 2889
$ws.Cells.Item(6, 16).Value = 963

# Row 7
$ws.Cells.Item(7, 4).Value = 44291
$ws.Cells.Item(7, 10).Value = 45
$ws.Cells.Item(7, 11).Value = 3000
$ws.Cells.Item(7, 12).Value = 3000
$ws.Cells.Item(7, 13).Value = 3000
$ws.Cells.Item(7, 16).Value = 1000

# Row 8
$ws.Cells.Item(8, 4).Value = 44193
$ws.Cells.Item(8, 10).Value = 70
$ws.Cells.Item(8, 11).Value = 3000
$ws.Cells.Item(8, 12).Value = 3000
$ws.Cells.Item(8, 13).Value = 3000
$ws.Cells.Item(8, 16).Value = 1000

# Row 9
$ws.Cells.Item(9, 4).Value = 44225
$ws.Cells.Item(9, 10).Value = 56
$ws.Cells.Item(9, 11).Value = 3000
$ws.Cells.Item(9, 12).Value = 3000
$ws.Cells.Item(9, 13).Value = 3000
$ws.Cells.Item(9, 16).Value = 1000

# Row 10
$ws.Cells.Item(10, 4).Value = 44179
$ws.Cells.Item(10, 10).Value = 78
$ws.Cells.Item(10, 11).Value = 3000
$ws.Cells.Item(10, 12).Value = 3000
$ws.Cells.Item(10, 13).Value = 3000
$ws.Cells.Item(10, 16).Value = 1000

# Row 11
$ws.Cells.Item(11, 4).Value = 44165
$ws.Cells.Item(11, 10).Value = 68
$ws.Cells.Item(11, 11).Value = 3000
$ws.Cells.Item(11, 12).Value = 3000
$ws.Cells.Item(11, 13).Value = 3000
$ws.Cells.Item(11, 16).Value = 1000

# Row 12 unchanged (44292, 40, 3000, 3000, 3000, 1000)

# Row 13
$ws.Cells.Item(13, 4).Value = 44243
$ws.Cells.Item(13, 10).Value = 45
$ws.Cells.Item(13, 11).Value = 3000
$ws.Cells.Item(13, 12).Value = 3000
$ws.Cells.Item(13, 13).Value = 3000
$ws.Cells.Item(13, 16).Value = 1000

# Row 14
$ws.Cells.Item(14, 4).Value = 44390
$ws.Cells.Item(14, 10).Value = 50
$ws.Cells.Item(14, 11).Value = 3000
$ws.Cells.Item(14, 12).Value = 3000
$ws.Cells.Item(14, 13).Value = 3000
$ws.Cells.Item(14, 16).Value = 1000

# Row 15
$ws.Cells.Item(15, 4).Value = 44221
$ws.Cells.Item(15, 10).Value = 50
$ws.Cells.Item(15, 11).Value = 2500
$ws.Cells.Item(15, 12).Value = 2500
$ws.Cells.Item(15, 13).Value = 2500
$ws.Cells.Item(15, 16).Value = 833

# Row 16
$ws.Cells.Item(16, 4).Value = 44242
$ws.Cells.Item(16, 10).Value = 95
$ws.Cells.Item(16, 11).Value = 2500
$ws.Cells.Item(16, 12).Value = 3000
$ws.Cells.Item(16, 13).Value = 2737
$ws.Cells.Item(16, 16).Value = 912

# Row 17
$ws.Cells.Item(17, 4).Value = 44166
$ws.Cells.Item(17, 10).Value = 45
$ws.Cells.Item(17, 11).Value = 2500
$ws.Cells.Item(17, 12).Value = 2500
$ws.Cells.Item(17, 13).Value = 2500
$ws.Cells.Item(17, 16).Value = 833

# Row 18
$ws.Cells.Item(18, 4).Value = 44340
$ws.Cells.Item(18, 10).Value = 54
$ws.Cells.Item(18, 11).Value = 3000
$ws.Cells.Item(18, 12).Value = 3000
$ws.Cells.Item(18, 13).Value = 3000
$ws.Cells.Item(18, 16).Value = 1000

# Row 19
$ws.Cells.Item(19, 4).Value = 44224
$ws.Cells.Item(19, 10).Value = 67
$ws.Cells.Item(19, 11).Value = 3000
$ws.Cells.Item(19, 12).Value = 3000
$ws.Cells.Item(19, 13).Value = 3000
$ws.Cells.Item(19, 16).Value = 1000

Write-Host "Done updating rows."
